$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 48 (shifts U9, Y1, Y2 rows down by one)
$ws.Rows.Item(48).Insert()

# Fill in the new row 48 with the FRAM device sharing the U8 footprint
$ws.Range("A48").Value = "FM24CL04B"
$ws.Range("B48").Value = "U8"
$ws.Range("C48").Value = "jlcpcb_smt:SOIC-8_3.9x4.9mm_P1.27mm"
$ws.Range("E48").Value = "C66016"
$ws.Range("D48").Value = "428-3741-1-ND"

# Update the view state to match the recorded selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("D48").Select()
